$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H33").Value = 522.7222
$ws.Range("I33").Value = 504.3125
$ws.Range("K33").Value = 504.3125
$ws.Range("M33").Value = -275.3125
$ws.Range("H125").Value = 9473.416999999999
$ws.Range("I125").Value = 1242
$ws.Range("J125").Value = 15353
$ws.Range("K125").Value = 11178
$ws.Range("L125").Value = 138177
$ws.Range("M125").Value = -8718
$ws.Range("N125").Value = -143097
$ws.Range("H132").Value = 20629.963
$ws.Range("I132").Value = 1395.6111
$ws.Range("K132").Value = 4186.8333
$ws.Range("M132").Value = -1656.8333
$ws.Range("H137").Value = 1343.2046
$ws.Range("I137").Value = 874.3226
$ws.Range("K137").Value = 2622.9678
$ws.Range("M137").Value = -72.9677999999999
$ws.Range("H138").Value = 1749.1285
$ws.Range("I138").Value = 879.8
$ws.Range("J138").Value = 2618.457
$ws.Range("K138").Value = 2639.4
$ws.Range("L138").Value = 7855.370999999999
$ws.Range("M138").Value = 2500.6
$ws.Range("N138").Value = -18135.371

# --- Sheet: ARM ---
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 9721.034
$ws.Range("I32").Value = 8206.875
$ws.Range("K32").Value = 8206.875
$ws.Range("M32").Value = -7919.875
$ws.Range("H37").Value = 7000
$ws.Range("I37").Value = 7000
$ws.Range("K37").Value = 7000
$ws.Range("M37").Value = -6727
$ws.Range("H45").Value = 2888.4443
$ws.Range("I45").Value = 2008.3334
$ws.Range("J45").Value = 3328.5
$ws.Range("K45").Value = 2008.3334
$ws.Range("L45").Value = 3328.5
$ws.Range("M45").Value = -1631.3334
$ws.Range("N45").Value = -4082.5
$ws.Range("H61").Value = 3459.739
$ws.Range("I61").Value = 1699.6154
$ws.Range("K61").Value = 1699.6154
$ws.Range("M61").Value = -1487.6154
$ws.Range("H74").Value = 1231.5
$ws.Range("I74").Value = 1106.2142
$ws.Range("J74").Value = 1523.8334
$ws.Range("K74").Value = 1106.2142
$ws.Range("L74").Value = 1523.8334
$ws.Range("M74").Value = -232.2141999999999
$ws.Range("N74").Value = -3271.8334
$ws.Range("H77").Value = 1231.5
$ws.Range("I77").Value = 1106.2142
$ws.Range("J77").Value = 1523.8334
$ws.Range("K77").Value = 5531.071
$ws.Range("L77").Value = 7619.166999999999
$ws.Range("M77").Value = -1163.071
$ws.Range("N77").Value = -16355.167
$ws.Range("H97").Value = 1074.4736
$ws.Range("I97").Value = 882.25
$ws.Range("J97").Value = 2099.6667
$ws.Range("K97").Value = 882.25
$ws.Range("L97").Value = 2099.6667
$ws.Range("M97").Value = -386.25
$ws.Range("N97").Value = -3091.6667
$ws.Range("H102").Value = 2389.2
$ws.Range("I102").Value = 1486.625
$ws.Range("J102").Value = 5999.5
$ws.Range("K102").Value = 1486.625
$ws.Range("L102").Value = 5999.5
$ws.Range("M102").Value = 135.375
$ws.Range("N102").Value = -9243.5
$ws.Range("H104").Value = 59000
$ws.Range("J104").Value = 59000
$ws.Range("L104").Value = 59000
$ws.Range("N104").Value = -65988
$ws.Range("H135").Value = 57395.312
$ws.Range("J135").Value = 57395.312
$ws.Range("L135").Value = 57395.312
$ws.Range("N135").Value = -67535.31200000001
$ws.Range("H136").Value = 3459.739
$ws.Range("I136").Value = 1699.6154
$ws.Range("K136").Value = 5098.8462
$ws.Range("M136").Value = -2548.8462

# --- Sheet: BSM ---
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H86").Value = 3258.1072
$ws.Range("I86").Value = 1572.5
$ws.Range("J86").Value = 6292.2
$ws.Range("K86").Value = 1572.5
$ws.Range("L86").Value = 6292.2
$ws.Range("M86").Value = -449.5
$ws.Range("N86").Value = -8538.200000000001
$ws.Range("H89").Value = 3258.1072
$ws.Range("I89").Value = 1572.5
$ws.Range("J89").Value = 6292.2
$ws.Range("K89").Value = 7862.5
$ws.Range("L89").Value = 31461
$ws.Range("M89").Value = -2246.5
$ws.Range("N89").Value = -42693
$ws.Range("H94").Value = 4098.9414
$ws.Range("I94").Value = 2845.4666
$ws.Range("J94").Value = 13500
$ws.Range("K94").Value = 2845.4666
$ws.Range("L94").Value = 13500
$ws.Range("M94").Value = -2394.4666
$ws.Range("N94").Value = -14402
$ws.Range("H105").Value = 2975
$ws.Range("I105").Value = 2975
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2975
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1228
$ws.Range("N105").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H7").Value = 274
$ws.Range("H31").Value = 3421.261
$ws.Range("I31").Value = 1285.6666
$ws.Range("J31").Value = 4794.143
$ws.Range("K31").Value = 1285.6666
$ws.Range("L31").Value = 4794.143
$ws.Range("M31").Value = -990.6666
$ws.Range("N31").Value = -5384.143
$ws.Range("H34").Value = 3421.261
$ws.Range("I34").Value = 1285.6666
$ws.Range("J34").Value = 4794.143
$ws.Range("K34").Value = 1285.6666
$ws.Range("L34").Value = 4794.143
$ws.Range("M34").Value = -1083.6666
$ws.Range("N34").Value = -5198.143
$ws.Range("H58").Value = 1740.0769
$ws.Range("I58").Value = 922.2143
$ws.Range("J58").Value = 2694.25
$ws.Range("K58").Value = 922.2143
$ws.Range("L58").Value = 2694.25
$ws.Range("M58").Value = -719.2143
$ws.Range("N58").Value = -3100.25
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H120").Value = 604442.25
$ws.Range("J120").Value = 604442.25
$ws.Range("L120").Value = 604442.25
$ws.Range("N120").Value = -611700.25
$ws.Range("H134").Value = 2834.8196
$ws.Range("I134").Value = 1614.9286
$ws.Range("K134").Value = 4844.7858
$ws.Range("M134").Value = -2309.7858
$ws.Range("H136").Value = 1740.0769
$ws.Range("I136").Value = 922.2143
$ws.Range("J136").Value = 2694.25
$ws.Range("K136").Value = 2766.6429
$ws.Range("L136").Value = 8082.75
$ws.Range("M136").Value = -216.6428999999998
$ws.Range("N136").Value = -13182.75

# --- Sheet: CUL ---
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H75").Value = 2931.3076
$ws.Range("I75").Value = 1022
$ws.Range("J75").Value = 3278.4546
$ws.Range("K75").Value = 3066
$ws.Range("L75").Value = 9835.363799999999
$ws.Range("M75").Value = -2068
$ws.Range("N75").Value = -11831.3638
$ws.Range("H78").Value = 2931.3076
$ws.Range("I78").Value = 1022
$ws.Range("J78").Value = 3278.4546
$ws.Range("K78").Value = 9198
$ws.Range("L78").Value = 29506.0914
$ws.Range("M78").Value = -4206
$ws.Range("N78").Value = -39490.0914
$ws.Range("H80").Value = 4239.4
$ws.Range("J80").Value = 4538.222
$ws.Range("L80").Value = 13614.666
$ws.Range("N80").Value = -15486.666
$ws.Range("H83").Value = 4239.4
$ws.Range("J83").Value = 4538.222
$ws.Range("L83").Value = 40843.998
$ws.Range("N83").Value = -50203.998
$ws.Range("H115").Value = 3002.2666
$ws.Range("J115").Value = 4857.143
$ws.Range("L115").Value = 14571.429
$ws.Range("N115").Value = -16921.429
$ws.Range("H139").Value = 2995.5715
$ws.Range("I139").Value = 2995.5715
$ws.Range("K139").Value = 8986.7145
$ws.Range("M139").Value = -3846.7145

# --- Sheet: GSM ---
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H10").Value = 6667.6665
$ws.Range("I10").Value = 7501.5
$ws.Range("K10").Value = 7501.5
$ws.Range("M10").Value = -7332.5
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H14").Value = 1561666.6
$ws.Range("I14").Value = 3750000
$ws.Range("J14").Value = 467500
$ws.Range("K14").Value = 3750000
$ws.Range("L14").Value = 467500
$ws.Range("M14").Value = -3749832
$ws.Range("N14").Value = -467836

# --- Sheet: LTW ---
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H22").Value = 696
$ws.Range("I22").Value = 681.4286
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 681.4286
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -386.4286
$ws.Range("N22").Value = -1490
$ws.Range("H27").Value = 696
$ws.Range("I27").Value = 681.4286
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 681.4286
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -574.4286
$ws.Range("N27").Value = -1114
$ws.Range("H46").Value = 2033
$ws.Range("J46").Value = 1200
$ws.Range("L46").Value = 1200
$ws.Range("N46").Value = -1576
$ws.Range("H68").Value = 6832.6665
$ws.Range("I68").Value = 6119.6
$ws.Range("J68").Value = 7342
$ws.Range("K68").Value = 6119.6
$ws.Range("L68").Value = 7342
$ws.Range("M68").Value = -5370.6
$ws.Range("N68").Value = -8840
$ws.Range("H71").Value = 6832.6665
$ws.Range("I71").Value = 6119.6
$ws.Range("J71").Value = 7342
$ws.Range("K71").Value = 30598
$ws.Range("L71").Value = 36710
$ws.Range("M71").Value = -26854
$ws.Range("N71").Value = -44198

# --- Sheet: WVR ---
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H124").Value = 499950
$ws.Range("J124").Value = 499950
$ws.Range("L124").Value = 499950
$ws.Range("N124").Value = -509770
$ws.Range("H126").Value = 2647.7693
$ws.Range("I126").Value = 2549.5557
$ws.Range("J126").Value = 2868.75
$ws.Range("K126").Value = 7648.6671
$ws.Range("L126").Value = 8606.25
$ws.Range("M126").Value = -5178.6671
$ws.Range("N126").Value = -13546.25
$ws.Range("H132").Value = 1412.238
$ws.Range("I132").Value = 950.0755
$ws.Range("K132").Value = 2850.2265
$ws.Range("M132").Value = -320.2265000000002
$ws.Range("H136").Value = 1672.5763
$ws.Range("I136").Value = 838.6739
$ws.Range("K136").Value = 2516.0217
$ws.Range("M136").Value = 33.97829999999976

